# Update Jogos_do_Dia_Betfair_Back_Lay_2026-01-08 with the latest odds feed.
# Rows 2-6 get refreshed odds; two new Algerian Ligue 1 fixtures are inserted
# (rows 3 and 7), which pushes the Italian Serie A / English Premier League
# rows further down the sheet (through row 9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Saudi Professional League | Al-Hilal vs Al-Hazm (KSA)
$ws.Range("A2").Value = "Saudi Professional League"
$ws.Range("B2").Value = "'2026-01-08"
$ws.Range("C2").Value = "11:55:00"
$ws.Range("D2").Value = "Al-Hilal"
$ws.Range("E2").Value = "Al-Hazm (KSA)"
$ws.Range("F2").Value = 1.17
$ws.Range("G2").Value = 1.2
$ws.Range("H2").Value = 15.5
$ws.Range("I2").Value = 22
$ws.Range("J2").Value = 9
$ws.Range("K2").Value = 11
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 8.4
$ws.Range("O2").Value = 1.1
$ws.Range("P2").Value = 3.75
$ws.Range("Q2").Value = 1.27
$ws.Range("R2").Value = 2.1
$ws.Range("S2").Value = 1.73
$ws.Range("T2").Value = 1.9
$ws.Range("U2").Value = 1.91
$ws.Range("V2").Value = 1.05
$ws.Range("W2").Value = 5.7
$ws.Range("X2").Value = 990
$ws.Range("Y2").Value = 990
$ws.Range("Z2").Value = 220
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 20
$ws.Range("AC2").Value = 30
$ws.Range("AD2").Value = 990
$ws.Range("AE2").Value = 260
$ws.Range("AF2").Value = 13.5
$ws.Range("AG2").Value = 16.5
$ws.Range("AH2").Value = 990
$ws.Range("AI2").Value = 170
$ws.Range("AJ2").Value = 11
$ws.Range("AK2").Value = 16.5
$ws.Range("AL2").Value = 980
$ws.Range("AM2").Value = 170
$ws.Range("AN2").Value = 2.76
$ws.Range("AO2").Value = 240

# Row 3: Algerian Ligue 1 | ES Setif vs Paradou
$ws.Range("A3").Value = "Algerian Ligue 1"
$ws.Range("B3").Value = "'2026-01-08"
$ws.Range("C3").Value = "13:45:00"
$ws.Range("D3").Value = "ES Setif"
$ws.Range("E3").Value = "Paradou"
$ws.Range("F3").Value = 1.04
$ws.Range("G3").Value = 1000
$ws.Range("H3").Value = 1.04
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 1.02
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 1.25
$ws.Range("O3").Value = 1.01
$ws.Range("P3").Value = 1.24
$ws.Range("Q3").Value = 1.01
$ws.Range("R3").Value = 1.12
$ws.Range("S3").Value = 1.05
$ws.Range("T3").Value = 1.04
$ws.Range("U3").Value = 1.04
$ws.Range("V3").Value = 1.01
$ws.Range("W3").Value = 1.01
$ws.Range("X3").Value = 990
$ws.Range("Y3").Value = 990
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 990
$ws.Range("AC3").Value = 990
$ws.Range("AD3").Value = 990
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 990
$ws.Range("AH3").Value = 990
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

# Row 4: Italian Serie A | US Cremonese vs Cagliari
$ws.Range("A4").Value = "Italian Serie A"
$ws.Range("B4").Value = "'2026-01-08"
$ws.Range("C4").Value = "14:30:00"
$ws.Range("D4").Value = "US Cremonese"
$ws.Range("E4").Value = "Cagliari"
$ws.Range("F4").Value = 2.54
$ws.Range("G4").Value = 2.56
$ws.Range("H4").Value = 3.35
$ws.Range("I4").Value = 3.45
$ws.Range("J4").Value = 3.15
$ws.Range("K4").Value = 3.2
$ws.Range("L4").Value = 1.52
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 3.05
$ws.Range("O4").Value = 1.46
$ws.Range("P4").Value = 1.67
$ws.Range("Q4").Value = 2.4
$ws.Range("R4").Value = 1.25
$ws.Range("S4").Value = 4.7
$ws.Range("T4").Value = 2
$ws.Range("U4").Value = 1.94
$ws.Range("V4").Value = 1.4
$ws.Range("W4").Value = 1.64
$ws.Range("X4").Value = 9.6
$ws.Range("Y4").Value = 10.5
$ws.Range("Z4").Value = 21
$ws.Range("AA4").Value = 60
$ws.Range("AB4").Value = 8.8
$ws.Range("AC4").Value = 7
$ws.Range("AD4").Value = 14.5
$ws.Range("AE4").Value = 46
$ws.Range("AF4").Value = 14.5
$ws.Range("AG4").Value = 12
$ws.Range("AH4").Value = 21
$ws.Range("AI4").Value = 65
$ws.Range("AJ4").Value = 36
$ws.Range("AK4").Value = 32
$ws.Range("AL4").Value = 55
$ws.Range("AM4").Value = 160
$ws.Range("AN4").Value = 32
$ws.Range("AO4").Value = 55

# Row 5: Saudi Professional League | Al Nassr vs Al-Quadisiya (KSA)
$ws.Range("A5").Value = "Saudi Professional League"
$ws.Range("B5").Value = "'2026-01-08"
$ws.Range("C5").Value = "14:30:00"
$ws.Range("D5").Value = "Al Nassr"
$ws.Range("E5").Value = "Al-Quadisiya (KSA)"
$ws.Range("F5").Value = 1.65
$ws.Range("G5").Value = 1.76
$ws.Range("H5").Value = 4.3
$ws.Range("I5").Value = 5.1
$ws.Range("J5").Value = 4.6
$ws.Range("K5").Value = 5.5
$ws.Range("L5").Value = 1.01
$ws.Range("M5").Value = 1.02
$ws.Range("N5").Value = 7.4
$ws.Range("O5").Value = 1.12
$ws.Range("P5").Value = 3.2
$ws.Range("Q5").Value = 1.36
$ws.Range("R5").Value = 1.76
$ws.Range("S5").Value = 1.81
$ws.Range("T5").Value = 1.46
$ws.Range("U5").Value = 2.48
$ws.Range("V5").Value = 1.24
$ws.Range("W5").Value = 2.28
$ws.Range("X5").Value = 42
$ws.Range("Y5").Value = 34
$ws.Range("Z5").Value = 50
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 18.5
$ws.Range("AC5").Value = 14
$ws.Range("AD5").Value = 22
$ws.Range("AE5").Value = 48
$ws.Range("AF5").Value = 17
$ws.Range("AG5").Value = 12
$ws.Range("AH5").Value = 17
$ws.Range("AI5").Value = 42
$ws.Range("AJ5").Value = 21
$ws.Range("AK5").Value = 16.5
$ws.Range("AL5").Value = 23
$ws.Range("AM5").Value = 60
$ws.Range("AN5").Value = 5.5
$ws.Range("AO5").Value = 28

# Row 6: Saudi Professional League | Al Najma Club vs Al-Ettifaq
$ws.Range("A6").Value = "Saudi Professional League"
$ws.Range("B6").Value = "'2026-01-08"
$ws.Range("C6").Value = "14:30:00"
$ws.Range("D6").Value = "Al Najma Club"
$ws.Range("E6").Value = "Al-Ettifaq"
$ws.Range("F6").Value = 4.4
$ws.Range("G6").Value = 5.3
$ws.Range("H6").Value = 1.79
$ws.Range("I6").Value = 1.97
$ws.Range("J6").Value = 3.65
$ws.Range("K6").Value = 4.1
$ws.Range("L6").Value = 1.01
$ws.Range("M6").Value = 1.07
$ws.Range("N6").Value = 3.45
$ws.Range("O6").Value = 1.32
$ws.Range("P6").Value = 1.84
$ws.Range("Q6").Value = 1.94
$ws.Range("R6").Value = 1.32
$ws.Range("S6").Value = 3.4
$ws.Range("T6").Value = 1.71
$ws.Range("U6").Value = 1.83
$ws.Range("V6").Value = 2.02
$ws.Range("W6").Value = 1.23
$ws.Range("X6").Value = 17
$ws.Range("Y6").Value = 10.5
$ws.Range("Z6").Value = 14
$ws.Range("AA6").Value = 26
$ws.Range("AB6").Value = 19.5
$ws.Range("AC6").Value = 10.5
$ws.Range("AD6").Value = 11
$ws.Range("AE6").Value = 22
$ws.Range("AF6").Value = 46
$ws.Range("AG6").Value = 24
$ws.Range("AH6").Value = 24
$ws.Range("AI6").Value = 46
$ws.Range("AJ6").Value = 140
$ws.Range("AK6").Value = 80
$ws.Range("AL6").Value = 85
$ws.Range("AM6").Value = 140
$ws.Range("AN6").Value = 95
$ws.Range("AO6").Value = 17

# Row 7: Algerian Ligue 1 | USM Alger vs MC Oran
$ws.Range("A7").Value = "Algerian Ligue 1"
$ws.Range("B7").Value = "'2026-01-08"
$ws.Range("C7").Value = "16:00:00"
$ws.Range("D7").Value = "USM Alger"
$ws.Range("E7").Value = "MC Oran"
$ws.Range("F7").Value = 1.04
$ws.Range("G7").Value = 1000
$ws.Range("H7").Value = 1.04
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 1.02
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 1.07
$ws.Range("Q7").Value = 1.01
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 0
$ws.Range("V7").Value = 0
$ws.Range("W7").Value = 0
$ws.Range("X7").Value = 0
$ws.Range("Y7").Value = 0
$ws.Range("Z7").Value = 0
$ws.Range("AA7").Value = 0
$ws.Range("AB7").Value = 0
$ws.Range("AC7").Value = 0
$ws.Range("AD7").Value = 0
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 0
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AI7").Value = 0
$ws.Range("AJ7").Value = 0
$ws.Range("AK7").Value = 0
$ws.Range("AL7").Value = 0
$ws.Range("AM7").Value = 0
$ws.Range("AN7").Value = 0
$ws.Range("AO7").Value = 0

# Row 8: Italian Serie A | AC Milan vs Genoa
$ws.Range("A8").Value = "Italian Serie A"
$ws.Range("B8").Value = "'2026-01-08"
$ws.Range("C8").Value = "16:45:00"
$ws.Range("D8").Value = "AC Milan"
$ws.Range("E8").Value = "Genoa"
$ws.Range("F8").Value = 1.45
$ws.Range("G8").Value = 1.46
$ws.Range("H8").Value = 9.6
$ws.Range("I8").Value = 10.5
$ws.Range("J8").Value = 4.7
$ws.Range("K8").Value = 4.8
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 3.55
$ws.Range("O8").Value = 1.37
$ws.Range("P8").Value = 1.85
$ws.Range("Q8").Value = 2.12
$ws.Range("R8").Value = 1.32
$ws.Range("S8").Value = 3.9
$ws.Range("T8").Value = 2.42
$ws.Range("U8").Value = 1.67
$ws.Range("V8").Value = 0
$ws.Range("W8").Value = 0
$ws.Range("X8").Value = 13.5
$ws.Range("Y8").Value = 25
$ws.Range("Z8").Value = 95
$ws.Range("AA8").Value = 590
$ws.Range("AB8").Value = 6.6
$ws.Range("AC8").Value = 10.5
$ws.Range("AD8").Value = 40
$ws.Range("AE8").Value = 1000
$ws.Range("AF8").Value = 7.2
$ws.Range("AG8").Value = 11
$ws.Range("AH8").Value = 34
$ws.Range("AI8").Value = 230
$ws.Range("AJ8").Value = 12
$ws.Range("AK8").Value = 18
$ws.Range("AL8").Value = 55
$ws.Range("AM8").Value = 320
$ws.Range("AN8").Value = 8.6
$ws.Range("AO8").Value = 1000

# Row 9: English Premier League | Arsenal vs Liverpool
$ws.Range("A9").Value = "English Premier League"
$ws.Range("B9").Value = "'2026-01-08"
$ws.Range("C9").Value = "17:00:00"
$ws.Range("D9").Value = "Arsenal"
$ws.Range("E9").Value = "Liverpool"
$ws.Range("F9").Value = 1.6
$ws.Range("G9").Value = 1.61
$ws.Range("H9").Value = 6.4
$ws.Range("I9").Value = 6.6
$ws.Range("J9").Value = 4.5
$ws.Range("K9").Value = 4.6
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 1.05
$ws.Range("N9").Value = 4.4
$ws.Range("O9").Value = 1.27
$ws.Range("P9").Value = 2.16
$ws.Range("Q9").Value = 1.83
$ws.Range("R9").Value = 1.46
$ws.Range("S9").Value = 3.05
$ws.Range("T9").Value = 1.89
$ws.Range("U9").Value = 2.06
$ws.Range("V9").Value = 0
$ws.Range("W9").Value = 0
$ws.Range("X9").Value = 18.5
$ws.Range("Y9").Value = 22
$ws.Range("Z9").Value = 55
$ws.Range("AA9").Value = 170
$ws.Range("AB9").Value = 8.6
$ws.Range("AC9").Value = 10
$ws.Range("AD9").Value = 24
$ws.Range("AE9").Value = 85
$ws.Range("AF9").Value = 9.2
$ws.Range("AG9").Value = 9.8
$ws.Range("AH9").Value = 21
$ws.Range("AI9").Value = 80
$ws.Range("AJ9").Value = 15
$ws.Range("AK9").Value = 15.5
$ws.Range("AL9").Value = 34
$ws.Range("AM9").Value = 110
$ws.Range("AN9").Value = 8.2
$ws.Range("AO9").Value = 100
